$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.493
$ws.Range("D5").Value = 0.591
$ws.Range("E5").Value = 0.611
$ws.Range("F5").Value = 0.681
$ws.Range("G5").Value = 0.681
$ws.Range("H5").Value = 0.6919999999999999

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.493
$ws.Range("F7").Value = 0.681
$ws.Range("G7").Value = 0.681
$ws.Range("H7").Value = 0.6919999999999999

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.488
$ws.Range("D8").Value = 0.645
$ws.Range("E8").Value = 0.678
$ws.Range("F8").Value = 0.697
$ws.Range("G8").Value = 0.727
$ws.Range("H8").Value = 0.744

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.291
$ws.Range("D9").Value = 0.62
$ws.Range("E9").Value = 0.654
$ws.Range("F9").Value = 0.644
$ws.Range("G9").Value = 0.678
$ws.Range("H9").Value = 0.6919999999999999
